$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch up formatting on row 12 (the last appointment) so it lines up with the
# rest of the table, and update its status to Declined.

# A/B/D/G already match the look used by the rows above (row 11 has the same
# per-column formatting), so just copy that across.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)

# C/E/F keep their own number formats/alignment but refresh the font.
$ws.Range("E12:F12").HorizontalAlignment = -4131
$ws.Range("C12").Font.Name = "Calibri"
$ws.Range("E12:F12").Font.Name = "Calibri"

$ws.Rows.Item(12).RowHeight = 17.25

$ws.Range("G12").Value = "Declined"
